# Update Betfair Back/Lay odds for 2026-01-05 games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "L2";  Value = 1.35 },
    @{ Cell = "N2";  Value = 4.3 },
    @{ Cell = "O2";  Value = 1.27 },
    @{ Cell = "P2";  Value = 2.14 },
    @{ Cell = "Q2";  Value = 1.79 },
    @{ Cell = "R2";  Value = 1.46 },
    @{ Cell = "S2";  Value = 2.96 },
    @{ Cell = "T2";  Value = 1.67 },
    @{ Cell = "U2";  Value = 2.32 },
    @{ Cell = "X2";  Value = 17 },
    @{ Cell = "Y2";  Value = 12 },
    @{ Cell = "AA2"; Value = 29 },
    @{ Cell = "AO2"; Value = 15.5 },

    @{ Cell = "Q3";  Value = 1.82 },

    @{ Cell = "G4";  Value = 1.51 },
    @{ Cell = "J4";  Value = 5.2 },

    @{ Cell = "G5";  Value = 1.83 },
    @{ Cell = "H5";  Value = 4.5 },

    @{ Cell = "P6";  Value = 1.7 },

    @{ Cell = "F7";  Value = 1.81 },
    @{ Cell = "G7";  Value = 2.08 },
    @{ Cell = "H7";  Value = 4.9 },
    @{ Cell = "I7";  Value = 6.8 },
    @{ Cell = "J7";  Value = 2.98 },
    @{ Cell = "K7";  Value = 3.55 },
    @{ Cell = "P7";  Value = 1.58 },
    @{ Cell = "Q7";  Value = 2.38 },

    @{ Cell = "P8";  Value = 1.73 },
    @{ Cell = "Q8";  Value = 2.12 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
